$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 8: Id=7, Name="Направлены на комиссию", Options="|sentOnCommission|"
# (matches the existing schema: Id, Name, Options)
#
# Column A holds a numeric-looking id ("7") that must be stored as TEXT
# (shared string), exactly like the existing rows (A2:A7 hold "1".."6" as
# text). A direct `.Value = "7"` assignment gets auto-coerced to a number
# by this host, so instead we render it through TEXT() in a scratch cell,
# copy/paste-special the *value* (not the formula) into place, which keeps
# it as text without touching the cell's style, then wipe the scratch cell.
$scratch = $ws.Cells.Item(10, 1)
$scratch.Formula = '=TEXT(7,"0")'
$scratch.Copy()
$ws.Cells.Item(8, 1).PasteSpecial(-4163)
$scratch.ClearContents()

$ws.Cells.Item(8, 2).Value = "Направлены на комиссию"
$ws.Cells.Item(8, 3).Value = "|sentOnCommission|"
